# Update the "取得日時" (retrieved datetime) timestamp in column A for all
# data rows (2-14) on the "ランサーズ" sheet to reflect the new scrape run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-24 18:32:15"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
